$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3926.625
$ws.Range("I74").Value = 3884.7
$ws.Range("J74").Value = 3996.5
$ws.Range("K74").Value = 3884.7
$ws.Range("L74").Value = 3996.5
$ws.Range("M74").Value = -2948.7
$ws.Range("N74").Value = -5868.5
$ws.Range("H77").Value = 3926.625
$ws.Range("I77").Value = 3884.7
$ws.Range("J77").Value = 3996.5
$ws.Range("K77").Value = 19423.5
$ws.Range("L77").Value = 19982.5
$ws.Range("M77").Value = -14743.5
$ws.Range("N77").Value = -29342.5
$ws.Range("H88").Value = 5372.1113
$ws.Range("J88").Value = 6081.2856
$ws.Range("L88").Value = 6081.2856
$ws.Range("N88").Value = -6893.2856
$ws.Range("H91").Value = 5372.1113
$ws.Range("J91").Value = 6081.2856
$ws.Range("L91").Value = 6081.2856
$ws.Range("N91").Value = -8889.285599999999
$ws.Range("H92").Value = 59773.234
$ws.Range("I92").Value = 125912.625
$ws.Range("J92").Value = 982.6667
$ws.Range("K92").Value = 125912.625
$ws.Range("L92").Value = 982.6667
$ws.Range("M92").Value = -124664.625
$ws.Range("N92").Value = -3478.6667
$ws.Range("H96").Value = 1559.5714
$ws.Range("I96").Value = 605.6667
$ws.Range("J96").Value = 2275
$ws.Range("K96").Value = 1817.0001
$ws.Range("L96").Value = 6825
$ws.Range("M96").Value = -444.0001
$ws.Range("N96").Value = -9571
$ws.Range("H99").Value = 2322.2727
$ws.Range("J99").Value = 6019.25
$ws.Range("L99").Value = 18057.75
$ws.Range("N99").Value = -21053.75
$ws.Range("H100").Value = 2826.76
$ws.Range("I100").Value = 2076.5
$ws.Range("J100").Value = 3519.3076
$ws.Range("K100").Value = 2076.5
$ws.Range("L100").Value = 3519.3076
$ws.Range("M100").Value = -1535.5
$ws.Range("N100").Value = -4601.3076
$ws.Range("H101").Value = 586
$ws.Range("I101").Value = 496.42856
$ws.Range("J101").Value = 899.5
$ws.Range("K101").Value = 1489.28568
$ws.Range("L101").Value = 2698.5
$ws.Range("M101").Value = 132.71432
$ws.Range("N101").Value = -5942.5
$ws.Range("H104").Value = 97
$ws.Range("I104").Value = 97
$ws.Range("K104").Value = 291
$ws.Range("M104").Value = 1456
$ws.Range("H127").Value = 1766.5
$ws.Range("I127").Value = 1766.5
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 5299.5
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -339.5
$ws.Range("N127").Value = ""
$ws.Range("H129").Value = 1250
$ws.Range("J129").Value = 1250
$ws.Range("L129").Value = 3750
$ws.Range("N129").Value = -13750
$ws.Range("H132").Value = 3384.8286
$ws.Range("I132").Value = 3384.8286
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10154.4858
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7624.485799999999
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19851.451
$ws.Range("J32").Value = 12960.8
$ws.Range("L32").Value = 12960.8
$ws.Range("N32").Value = -13534.8
$ws.Range("H97").Value = 5648.185
$ws.Range("I97").Value = 5950.25
$ws.Range("K97").Value = 5950.25
$ws.Range("M97").Value = -5454.25
$ws.Range("H102").Value = 2791.2778
$ws.Range("I102").Value = 2425.4119
$ws.Range("K102").Value = 2425.4119
$ws.Range("M102").Value = -803.4119000000001
$ws.Range("H122").Value = 1507.6316
$ws.Range("I122").Value = 973.0769
$ws.Range("J122").Value = 2665.8333
$ws.Range("K122").Value = 2919.2307
$ws.Range("L122").Value = 7997.499899999999
$ws.Range("M122").Value = -469.2307000000001
$ws.Range("N122").Value = -12897.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3265.0588
$ws.Range("I105").Value = 3233.0667
$ws.Range("K105").Value = 3233.0667
$ws.Range("M105").Value = -1486.0667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3418.1667
$ws.Range("J99").Value = 3435.8572
$ws.Range("L99").Value = 3435.8572
$ws.Range("N99").Value = -6431.8572
$ws.Range("H126").Value = 3418.1667
$ws.Range("J126").Value = 3435.8572
$ws.Range("L126").Value = 10307.5716
$ws.Range("N126").Value = -15247.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1472.5834
$ws.Range("I139").Value = 1151.909
$ws.Range("K139").Value = 3455.727
$ws.Range("M139").Value = 1684.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 784.7857
$ws.Range("I97").Value = 597.7
$ws.Range("J97").Value = 1252.5
$ws.Range("K97").Value = 597.7
$ws.Range("L97").Value = 1252.5
$ws.Range("M97").Value = -101.7
$ws.Range("N97").Value = -2244.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 44556.117
$ws.Range("I22").Value = 69822.875
$ws.Range("J22").Value = 4129.3
$ws.Range("K22").Value = 69822.875
$ws.Range("L22").Value = 4129.3
$ws.Range("M22").Value = -69527.875
$ws.Range("N22").Value = -4719.3
$ws.Range("H27").Value = 44556.117
$ws.Range("I27").Value = 69822.875
$ws.Range("J27").Value = 4129.3
$ws.Range("K27").Value = 69822.875
$ws.Range("L27").Value = 4129.3
$ws.Range("M27").Value = -69715.875
$ws.Range("N27").Value = -4343.3
$ws.Range("H93").Value = 2013.1842
$ws.Range("I93").Value = 1757.4117
$ws.Range("J93").Value = 2220.238
$ws.Range("K93").Value = 1757.4117
$ws.Range("L93").Value = 2220.238
$ws.Range("M93").Value = -509.4117000000001
$ws.Range("N93").Value = -4716.237999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 113245.1
$ws.Range("I62").Value = 9660.666999999999
$ws.Range("J62").Value = 157638.42
$ws.Range("K62").Value = 9660.666999999999
$ws.Range("L62").Value = 157638.42
$ws.Range("M62").Value = -9036.666999999999
$ws.Range("N62").Value = -158886.42
$ws.Range("H65").Value = 113245.1
$ws.Range("I65").Value = 9660.666999999999
$ws.Range("J65").Value = 157638.42
$ws.Range("K65").Value = 48303.335
$ws.Range("L65").Value = 788192.1000000001
$ws.Range("M65").Value = -45183.335
$ws.Range("N65").Value = -794432.1000000001
$ws.Range("H100").Value = 687.9231
$ws.Range("I100").Value = 604.5454999999999
$ws.Range("K100").Value = 1209.091
$ws.Range("M100").Value = -668.0909999999999
$ws.Range("H132").Value = 56912.688
$ws.Range("I132").Value = 61892.03
$ws.Range("K132").Value = 185676.09
$ws.Range("M132").Value = -183146.09
